$d = $word.ActiveDocument

# --- 1. Title paragraph "Team Member Profiles": drop bold, bump size 24 -> 36 (half-points: 12pt -> 18pt) ---
$title = $d.Paragraphs(1)
$title.Range.Font.Size = 18
$title.Range.Font.SizeBi = 18
$title.Range.Font.Bold = 0
$title.Range.Font.BoldBi = 0

# --- 2. Add an extra blank paragraph before "Name: Daniel Butler" (after the Xavier Ruyle block) ---
# Paragraph 6 is the second of the two existing blank paragraphs that sit right before
# the "Name: Daniel Butler" paragraph.
$pDaniel = $d.Paragraphs(6)
$pDaniel.Range.InsertParagraphAfter()

# --- 3. Add an extra blank paragraph right after "Major: Computer Science" (Michael Hoopes block),
#        before the existing blank paragraph that precedes "Name: Connor Williamson" ---
$pMichaelMajor = $d.Paragraphs(13)
$pMichaelMajor.Range.InsertParagraphAfter()
